$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the baud rate from 6 MBaud to 8 MBaud
$ws.Range("B7").Formula = "=8*10^6"

# Update the sheet view: scroll back to top and change selection to B8
$ws.Range("B8").Select()
